# "easy stages enemy added"
# Appends the roster of enemies for the Easy stage (Fruit Village / Veggie
# Forest area) below the existing table, in column A starting at row 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$enemies = @(
    "Mole",
    "Mice",
    "Prideful Worm",
    "Beetle 1",
    "Beetle 2",
    "Beetle 3",
    "Queen Bitter Gourd",
    "Insect Queen"
)

$startRow = 19
for ($i = 0; $i -lt $enemies.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $enemies[$i]
}

# The boss entry (row 22, "Beetle 1") gets its own left-aligned, non-bold
# Arial style instead of the bold/centered style that used to live on A22.
$bossCell = $ws.Range("A22")
$bossCell.Font.Name = "Arial"
$bossCell.Font.Bold = $false
$bossCell.HorizontalAlignment = -4131  # xlLeft

# Rows 20-22 keep the slightly taller 17.25pt row height.
$ws.Rows.Item(20).RowHeight = 17.25
$ws.Rows.Item(21).RowHeight = 17.25
$ws.Rows.Item(22).RowHeight = 17.25

# Restore the view/selection state recorded after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A22").Select()
